# Refresh crypto price/volume table (rows 2-51) with latest scrape values.
# Mirrors GitHub Actions commit "Updated cryptos list on Mon Dec 11 15:25:55 UTC 2023".
#
# Numeric-looking Price values (e.g. "244.13") are written with a leading
# backtick-escaped apostrophe -- Excel's standard "force text" quote prefix --
# so they stay text cells (matching the source file's inlineStr cells) instead
# of being auto-coerced to numbers. The apostrophe itself is not stored in the
# cell value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.057.26"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "2.227.09"
$ws.Range("E3").Value = "  -5.36%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'244.13"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "'0.627"
$ws.Range("E6").Value = "  -6.32%  "
$ws.Range("D7").Value = "'68.50"
$ws.Range("E7").Value = "  -8.07%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  -7.83%  "
$ws.Range("D10").Value = "'0.0960"
$ws.Range("E10").Value = "  -4.96%  "
$ws.Range("D11").Value = "'58.15"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("D12").Value = "'35.57"
$ws.Range("E12").Value = "  +6.56%  "
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "'6.73"
$ws.Range("E14").Value = "  -7.75%  "
$ws.Range("D15").Value = "2.555.69"
$ws.Range("E15").Value = "  -5.44%  "
$ws.Range("D16").Value = "'14.79"
$ws.Range("E16").Value = "  -8.81%  "
$ws.Range("D17").Value = "'0.850"
$ws.Range("E17").Value = "  -6.32%  "
$ws.Range("D18").Value = "2.224.12"
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("D19").Value = "41.962.21"
$ws.Range("E19").Value = "  -4.40%  "
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("E20").Value = "  -7.57%  "
$ws.Range("D21").Value = "'72.56"
$ws.Range("E21").Value = "  -7.32%  "
$ws.Range("D22").Value = "'6.13"
$ws.Range("E22").Value = "  -7.99%  "
$ws.Range("D23").Value = "'234.71"
$ws.Range("E23").Value = "  -7.36%  "
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  +10.48%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "'3.62"
$ws.Range("E26").Value = "  -4.41%  "
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("D28").Value = "'2.24"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").Value = "'9.89"
$ws.Range("E29").Value = "  -5.72%  "
$ws.Range("D30").Value = "'171.38"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").Value = "'20.44"
$ws.Range("E31").Value = "  -8.46%  "
$ws.Range("D32").Value = "'0.121"
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("E33").Value = "  -6.64%  "
$ws.Range("D34").Value = "'0.0713"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").Value = "'5.20"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "'4.66"
$ws.Range("E36").Value = "  -8.48%  "
$ws.Range("D37").Value = "'3.89"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "'22.50"
$ws.Range("E38").Value = "  +18.21%  "
$ws.Range("D39").Value = "'0.0279"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("D40").Value = "'2.30"
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("D41").Value = "'66.57"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("D42").Value = "'5.82"
$ws.Range("E42").Value = "  -9.44%  "
$ws.Range("D43").Value = "'4.96"
$ws.Range("E43").Value = "  -11.15%  "
$ws.Range("D44").Value = "'9.02"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "'4.55"
$ws.Range("E48").Value = "  +7.30%  "
$ws.Range("D49").Value = "'1.19"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").Value = "'10.14"
$ws.Range("E50").Value = "  +7.51%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.11"
$ws.Range("E51").Value = "  -4.21%  "
